$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 32.21267
$ws.Cells.Item(2,8).Value = 96.63801000000001
$ws.Cells.Item(2,9).Value = 0.7096649552378644
$ws.Cells.Item(2,10).Value = 0.7096649552378644
$ws.Cells.Item(2,13).Value = 0.379829
$ws.Cells.Item(2,14).Value = 1.139487
$ws.Cells.Item(2,15).Value = 0.01205513428578339
$ws.Cells.Item(2,16).Value = 0.01205513428578339
$ws.Cells.Item(2,17).Value = 12.23530623343
$ws.Cells.Item(2,18).Value = 110.11775610087
$ws.Cells.Item(2,19).Value = 0.008555106333306911
$ws.Cells.Item(2,20).Value = 0.008555106333306911
$ws.Cells.Item(3,7).Value = 32.21267
$ws.Cells.Item(3,8).Value = 96.63801000000001
$ws.Cells.Item(3,9).Value = 0.7096649552378644
$ws.Cells.Item(3,10).Value = 0.7096649552378644
$ws.Cells.Item(3,13).Value = 5.511159333333334
$ws.Cells.Item(3,14).Value = 16.533478
$ws.Cells.Item(3,15).Value = 0.17491493760003
$ws.Cells.Item(3,16).Value = 0.17491493760003
$ws.Cells.Item(3,17).Value = 177.5291569220867
$ws.Cells.Item(3,18).Value = 1597.76241229878
$ws.Cells.Item(3,19).Value = 0.1241310013623591
$ws.Cells.Item(3,20).Value = 0.1241310013623591
$ws.Cells.Item(4,7).Value = 32.21267
$ws.Cells.Item(4,8).Value = 96.63801000000001
$ws.Cells.Item(4,9).Value = 0.7096649552378644
$ws.Cells.Item(4,10).Value = 0.7096649552378644
$ws.Cells.Item(4,13).Value = 1.075784333333333
$ws.Cells.Item(4,14).Value = 3.227353
$ws.Cells.Item(4,15).Value = 0.034143587248144
$ws.Cells.Item(4,16).Value = 0.034143587248144
$ws.Cells.Item(4,17).Value = 34.65388572083667
$ws.Cells.Item(4,18).Value = 311.88497148753
$ws.Cells.Item(4,19).Value = 0.02423050731611423
$ws.Cells.Item(4,20).Value = 0.02423050731611423
$ws.Cells.Item(5,7).Value = 32.21267
$ws.Cells.Item(5,8).Value = 96.63801000000001
$ws.Cells.Item(5,9).Value = 0.7096649552378644
$ws.Cells.Item(5,10).Value = 0.7096649552378644
$ws.Cells.Item(5,13).Value = 24.54088133333333
$ws.Cells.Item(5,14).Value = 73.622644
$ws.Cells.Item(5,15).Value = 0.7788863408660427
$ws.Cells.Item(5,16).Value = 0.7788863408660427
$ws.Cells.Item(5,17).Value = 790.5273118998267
$ws.Cells.Item(5,18).Value = 7114.74580709844
$ws.Cells.Item(5,19).Value = 0.5527483402260842
$ws.Cells.Item(5,20).Value = 0.5527483402260842
$ws.Cells.Item(6,9).Value = 0.2527239295880077
$ws.Cells.Item(6,10).Value = 0.2527239295880077
$ws.Cells.Item(6,13).Value = 0.379829
$ws.Cells.Item(6,14).Value = 1.139487
$ws.Cells.Item(6,15).Value = 0.01205513428578339
$ws.Cells.Item(6,16).Value = 0.01205513428578339
$ws.Cells.Item(6,17).Value = 4.357203562332667
$ws.Cells.Item(6,18).Value = 39.21483206099401
$ws.Cells.Item(6,19).Value = 0.003046620908414298
$ws.Cells.Item(6,20).Value = 0.003046620908414298
$ws.Cells.Item(7,9).Value = 0.2527239295880077
$ws.Cells.Item(7,10).Value = 0.2527239295880077
$ws.Cells.Item(7,13).Value = 5.511159333333334
$ws.Cells.Item(7,14).Value = 16.533478
$ws.Cells.Item(7,15).Value = 0.17491493760003
$ws.Cells.Item(7,16).Value = 0.17491493760003
$ws.Cells.Item(7,17).Value = 63.22119448431513
$ws.Cells.Item(7,18).Value = 568.9907503588361
$ws.Cells.Item(7,19).Value = 0.04420519037392073
$ws.Cells.Item(7,20).Value = 0.04420519037392073
$ws.Cells.Item(8,9).Value = 0.2527239295880077
$ws.Cells.Item(8,10).Value = 0.2527239295880077
$ws.Cells.Item(8,13).Value = 1.075784333333333
$ws.Cells.Item(8,14).Value = 3.227353
$ws.Cells.Item(8,15).Value = 0.034143587248144
$ws.Cells.Item(8,16).Value = 0.034143587248144
$ws.Cells.Item(8,17).Value = 12.34084635323178
$ws.Cells.Item(8,18).Value = 111.067617179086
$ws.Cells.Item(8,19).Value = 0.008628901539581943
$ws.Cells.Item(8,20).Value = 0.008628901539581943
$ws.Cells.Item(9,9).Value = 0.2527239295880077
$ws.Cells.Item(9,10).Value = 0.2527239295880077
$ws.Cells.Item(9,13).Value = 24.54088133333333
$ws.Cells.Item(9,14).Value = 73.622644
$ws.Cells.Item(9,15).Value = 0.7788863408660427
$ws.Cells.Item(9,16).Value = 0.7788863408660427
$ws.Cells.Item(9,17).Value = 281.5204093641698
$ws.Cells.Item(9,18).Value = 2533.683684277528
$ws.Cells.Item(9,19).Value = 0.1968432167660907
$ws.Cells.Item(9,20).Value = 0.1968432167660907
$ws.Cells.Item(10,7).Value = 1.279382333333333
$ws.Cells.Item(10,8).Value = 3.838147
$ws.Cells.Item(10,9).Value = 0.02818558059040478
$ws.Cells.Item(10,10).Value = 0.02818558059040478
$ws.Cells.Item(10,13).Value = 0.379829
$ws.Cells.Item(10,14).Value = 1.139487
$ws.Cells.Item(10,15).Value = 0.01205513428578339
$ws.Cells.Item(10,16).Value = 0.01205513428578339
$ws.Cells.Item(10,17).Value = 0.4859465122876667
$ws.Cells.Item(10,18).Value = 4.373518610589
$ws.Cells.Item(10,19).Value = 0.0003397809589400994
$ws.Cells.Item(10,20).Value = 0.0003397809589400994
$ws.Cells.Item(11,7).Value = 1.279382333333333
$ws.Cells.Item(11,8).Value = 3.838147
$ws.Cells.Item(11,9).Value = 0.02818558059040478
$ws.Cells.Item(11,10).Value = 0.02818558059040478
$ws.Cells.Item(11,13).Value = 5.511159333333334
$ws.Cells.Item(11,14).Value = 16.533478
$ws.Cells.Item(11,15).Value = 0.17491493760003
$ws.Cells.Item(11,16).Value = 0.17491493760003
$ws.Cells.Item(11,17).Value = 7.050879887251778
$ws.Cells.Item(11,18).Value = 63.45791898526601
$ws.Cells.Item(11,19).Value = 0.004930079070191268
$ws.Cells.Item(11,20).Value = 0.004930079070191268
$ws.Cells.Item(12,7).Value = 1.279382333333333
$ws.Cells.Item(12,8).Value = 3.838147
$ws.Cells.Item(12,9).Value = 0.02818558059040478
$ws.Cells.Item(12,10).Value = 0.02818558059040478
$ws.Cells.Item(12,13).Value = 1.075784333333333
$ws.Cells.Item(12,14).Value = 3.227353
$ws.Cells.Item(12,15).Value = 0.034143587248144
$ws.Cells.Item(12,16).Value = 0.034143587248144
$ws.Cells.Item(12,17).Value = 1.376339470543444
$ws.Cells.Item(12,18).Value = 12.387055234891
$ws.Cells.Item(12,19).Value = 0.0009623568300280798
$ws.Cells.Item(12,20).Value = 0.0009623568300280798
$ws.Cells.Item(13,7).Value = 1.279382333333333
$ws.Cells.Item(13,8).Value = 3.838147
$ws.Cells.Item(13,9).Value = 0.02818558059040478
$ws.Cells.Item(13,10).Value = 0.02818558059040478
$ws.Cells.Item(13,13).Value = 24.54088133333333
$ws.Cells.Item(13,14).Value = 73.622644
$ws.Cells.Item(13,15).Value = 0.7788863408660427
$ws.Cells.Item(13,16).Value = 0.7788863408660427
$ws.Cells.Item(13,17).Value = 31.39717002229644
$ws.Cells.Item(13,18).Value = 282.5745302006679
$ws.Cells.Item(13,19).Value = 0.02195336373124534
$ws.Cells.Item(13,20).Value = 0.02195336373124534
$ws.Cells.Item(14,7).Value = 0.4278379999999999
$ws.Cells.Item(14,8).Value = 1.283514
$ws.Cells.Item(14,9).Value = 0.009425534583723031
$ws.Cells.Item(14,10).Value = 0.009425534583723031
$ws.Cells.Item(14,13).Value = 0.379829
$ws.Cells.Item(14,14).Value = 1.139487
$ws.Cells.Item(14,15).Value = 0.01205513428578339
$ws.Cells.Item(14,16).Value = 0.01205513428578339
$ws.Cells.Item(14,17).Value = 0.162505279702
$ws.Cells.Item(14,18).Value = 1.462547517318
$ws.Cells.Item(14,19).Value = 0.0001136260851220765
$ws.Cells.Item(14,20).Value = 0.0001136260851220766
$ws.Cells.Item(15,7).Value = 0.4278379999999999
$ws.Cells.Item(15,8).Value = 1.283514
$ws.Cells.Item(15,9).Value = 0.009425534583723031
$ws.Cells.Item(15,10).Value = 0.009425534583723031
$ws.Cells.Item(15,13).Value = 5.511159333333334
$ws.Cells.Item(15,14).Value = 16.533478
$ws.Cells.Item(15,15).Value = 0.17491493760003
$ws.Cells.Item(15,16).Value = 0.17491493760003
$ws.Cells.Item(15,17).Value = 2.357883386854667
$ws.Cells.Item(15,18).Value = 21.220950481692
$ws.Cells.Item(15,19).Value = 0.001648666793558838
$ws.Cells.Item(15,20).Value = 0.001648666793558838
$ws.Cells.Item(16,7).Value = 0.4278379999999999
$ws.Cells.Item(16,8).Value = 1.283514
$ws.Cells.Item(16,9).Value = 0.009425534583723031
$ws.Cells.Item(16,10).Value = 0.009425534583723031
$ws.Cells.Item(16,13).Value = 1.075784333333333
$ws.Cells.Item(16,14).Value = 3.227353
$ws.Cells.Item(16,15).Value = 0.034143587248144
$ws.Cells.Item(16,16).Value = 0.034143587248144
$ws.Cells.Item(16,17).Value = 0.4602614176046665
$ws.Cells.Item(16,18).Value = 4.142352758442
$ws.Cells.Item(16,19).Value = 0.0003218215624197459
$ws.Cells.Item(16,20).Value = 0.0003218215624197459
$ws.Cells.Item(17,7).Value = 0.4278379999999999
$ws.Cells.Item(17,8).Value = 1.283514
$ws.Cells.Item(17,9).Value = 0.009425534583723031
$ws.Cells.Item(17,10).Value = 0.009425534583723031
$ws.Cells.Item(17,13).Value = 24.54088133333333
$ws.Cells.Item(17,14).Value = 73.622644
$ws.Cells.Item(17,15).Value = 0.7788863408660427
$ws.Cells.Item(17,16).Value = 0.7788863408660427
$ws.Cells.Item(17,17).Value = 10.49952158789066
$ws.Cells.Item(17,18).Value = 94.49569429101598
$ws.Cells.Item(17,19).Value = 0.007341420142622371
$ws.Cells.Item(17,20).Value = 0.007341420142622371
